# This script applies the Sep 10 2024 "cryptos" price/volume refresh:
#  - updates Price (D) and Volume(1h) (E) figures for most rows
#  - re-sorts a few coin pairs that swapped rank order, so their
#    Coin/Link/Price/Volume columns (B/C/D/E) are fully replaced
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores values as plain text (e.g. "517.08", "1.00").
# Force text formatting before writing the new figures so Excel does not
# silently reinterpret them as numbers.
$priceTextCells = @("D5", "D6", "D8", "D10", "D11", "D12", "D13", "D14", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $priceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '56.924.03'
$ws.Range("E2").Value = '  +3.83%  '

$ws.Range("D3").Value = '2.347.40'
$ws.Range("E3").Value = '  +2.69%  '

$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").Value = '517.08'
$ws.Range("E5").Value = '  +1.71%  '

$ws.Range("D6").Value = '133.65'
$ws.Range("E6").Value = '  +2.75%  '

$ws.Range("E7").Value = '  +0.40%  '

$ws.Range("D8").Value = '0.535'
$ws.Range("E8").Value = '  +1.12%  '

$ws.Range("D9").Value = '2.343.66'
$ws.Range("E9").Value = '  +1.60%  '

$ws.Range("D10").Value = '0.103'
$ws.Range("E10").Value = '  +6.78%  '

$ws.Range("D11").Value = '0.154'
$ws.Range("E11").Value = '  +0.31%  '

$ws.Range("D12").Value = '5.20'
$ws.Range("E12").Value = '  +5.99%  '

$ws.Range("D13").Value = '0.340'
$ws.Range("E13").Value = '  -1.97%  '

$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = '23.73'
$ws.Range("E14").Value = '  +1.68%  '

$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '2.761.22'
$ws.Range("E15").Value = '  +2.54%  '

$ws.Range("D16").Value = '56.887.56'
$ws.Range("E16").Value = '  +3.69%  '

$ws.Range("D17").Value = '0.0000134'
$ws.Range("E17").Value = '  +2.36%  '

$ws.Range("D18").Value = '2.344.39'
$ws.Range("E18").Value = '  +2.03%  '

$ws.Range("D19").Value = '10.43'
$ws.Range("E19").Value = '  +0.73%  '

$ws.Range("D20").Value = '4.26'
$ws.Range("E20").Value = '  +2.23%  '

$ws.Range("D21").Value = '321.29'
$ws.Range("E21").Value = '  +4.49%  '

$ws.Range("D22").Value = '6.67'
$ws.Range("E22").Value = '  +4.69%  '

$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  +0.21%  '

$ws.Range("D24").Value = '60.80'
$ws.Range("E24").Value = '  +0.62%  '

$ws.Range("D25").Value = '0.996'
$ws.Range("E25").Value = '  +0.32%  '

$ws.Range("D26").Value = '0.159'
$ws.Range("E26").Value = '  +4.89%  '

$ws.Range("D27").Value = '7.76'
$ws.Range("E27").Value = '  +3.85%  '

$ws.Range("D28").Value = '171.31'
$ws.Range("E28").Value = '  -0.49%  '

$ws.Range("B29").Value = 'Fetch.AI'
$ws.Range("C29").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D29").Value = '1.21'
$ws.Range("E29").Value = '  +8.74%  '

$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '0.0₃0734'
$ws.Range("E30").Value = '  +2.72%  '

$ws.Range("D31").Value = '6.23'
$ws.Range("E31").Value = '  +2.61%  '

$ws.Range("D32").Value = '1.67'
$ws.Range("E32").Value = '  +2.22%  '

$ws.Range("D33").Value = '18.29'
$ws.Range("E33").Value = '  +1.61%  '

$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  +0.06%  '

$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.47%  '

$ws.Range("D36").Value = '0.954'
$ws.Range("E36").Value = '  +2.20%  '

$ws.Range("D37").Value = '1.25'
$ws.Range("E37").Value = '  +3.87%  '

$ws.Range("D38").Value = '3.99'
$ws.Range("E38").Value = '  +5.18%  '

$ws.Range("E39").Value = '  +6.85%  '

$ws.Range("D40").Value = '37.47'
$ws.Range("E40").Value = '  +2.28%  '

$ws.Range("D41").Value = '0.380'
$ws.Range("E41").Value = '  +0.45%  '

$ws.Range("D42").Value = '138.80'
$ws.Range("E42").Value = '  +10.61%  '

$ws.Range("D43").Value = '3.57'
$ws.Range("E43").Value = '  +4.39%  '

$ws.Range("D44").Value = '276.65'
$ws.Range("E44").Value = '  +9.54%  '

$ws.Range("D45").Value = '5.08'
$ws.Range("E45").Value = '  +3.81%  '

$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").Value = '0.0929'
$ws.Range("E46").Value = '  +2.49%  '

$ws.Range("B47").Value = 'Hedera'
$ws.Range("C47").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D47").Value = '0.0507'
$ws.Range("E47").Value = '  +1.68%  '

$ws.Range("D48").Value = '0.558'
$ws.Range("E48").Value = '  +1.14%  '

$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").Value = '0.0215'
$ws.Range("E49").Value = '  +3.71%  '

$ws.Range("B50").Value = 'Polygon'
$ws.Range("C50").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D50").Value = '0.380'
$ws.Range("E50").Value = '  +0.69%  '

$ws.Range("D51").Value = '16.82'
$ws.Range("E51").Value = '  +1.85%  '
